$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix A65: was stored as text, should become a number
$ws.Cells.Item(65, 1).Value = 51616191

# Add new row 66 for payment 71717172 (Cash) 2025-08-20T08:11:27
$a66 = $ws.Cells.Item(66, 1)
$a66.NumberFormat = "@"
$a66.Value = "71717172"
$a66.Style = "Normal"

$ws.Cells.Item(66, 2).Value = ""
$ws.Cells.Item(66, 3).Value = "Cash"
$ws.Cells.Item(66, 4).Value = "2025-08-20T08:11:27"
$ws.Cells.Item(66, 5).Value = 125
$ws.Cells.Item(66, 6).Value = ""
$ws.Cells.Item(66, 7).Value = 125
$ws.Cells.Item(66, 8).Value = 0
$ws.Cells.Item(66, 9).Value = 0
$ws.Cells.Item(66, 10).Value = 0
